$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L14").Value = 0
$ws.Range("N14").Value = -62.857142857142
$ws.Range("C15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 30
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = -14.285714285714
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -47.368421052631
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -22.222222222222
$ws.Range("F16").Value = 45
$ws.Range("G16").Value = 43
$ws.Range("H16").Value = 4.651162790697
$ws.Range("I16").Value = 424
$ws.Range("J16").Value = 435
$ws.Range("K16").Value = -2.528735632183
$ws.Range("L16").Value = 25.816023738872
$ws.Range("M16").Value = 10.704960835509
$ws.Range("N16").Value = -61.628959276018
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 19
$ws.Range("E17").Value = -36.842105263157
$ws.Range("F17").Value = 50
$ws.Range("G17").Value = 52
$ws.Range("H17").Value = -3.846153846153
$ws.Range("I17").Value = 689
$ws.Range("J17").Value = 610
$ws.Range("K17").Value = 12.950819672131
$ws.Range("L17").Value = 50.436681222707
$ws.Range("M17").Value = 94.084507042253
$ws.Range("N17").Value = 2.529761904761
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -7.142857142857
$ws.Range("I18").Value = 230
$ws.Range("J18").Value = 245
$ws.Range("K18").Value = -6.122448979591
$ws.Range("L18").Value = 9.004739336492
$ws.Range("M18").Value = -1.709401709401
$ws.Range("N18").Value = -81.466559226430
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 47
$ws.Range("H19").Value = 17.5
$ws.Range("I19").Value = 506
$ws.Range("J19").Value = 448
$ws.Range("K19").Value = 12.946428571428
$ws.Range("L19").Value = 14.739229024943
$ws.Range("M19").Value = 65.359477124183
$ws.Range("N19").Value = 15.525114155251
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -22.727272727272
$ws.Range("I20").Value = 311
$ws.Range("J20").Value = 279
$ws.Range("K20").Value = 11.469534050179
$ws.Range("L20").Value = 63.684210526315
$ws.Range("M20").Value = 177.678571428571
$ws.Range("N20").Value = -46.100519930675
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -11.363636363636
$ws.Range("F21").Value = 174
$ws.Range("G21").Value = 173
$ws.Range("H21").Value = 0.578034682080
$ws.Range("I21").Value = 2203
$ws.Range("J21").Value = 2053
$ws.Range("K21").Value = 7.306380905991
$ws.Range("L21").Value = 30.741839762611
$ws.Range("M21").Value = 55.799151343705
$ws.Range("N21").Value = -46.593939393939
$ws.Range("D22").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("L22").Value = -63.636363636363
$ws.Range("C23").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 22
$ws.Range("J23").Value = 28
$ws.Range("K23").Value = -21.428571428571
$ws.Range("L23").Value = 22.222222222222
$ws.Range("M23").Value = 37.5
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 38.095238095238
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = -5.617977528089
$ws.Range("I24").Value = 987
$ws.Range("J24").Value = 1162
$ws.Range("K24").Value = -15.060240963855
$ws.Range("L24").Value = 15.034965034965
$ws.Range("M24").Value = 6.587473002159
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 7.142857142857
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = 14.545454545454
$ws.Range("I25").Value = 897
$ws.Range("J25").Value = 839
$ws.Range("K25").Value = 6.912991656734
$ws.Range("L25").Value = 21.875
$ws.Range("M25").Value = -3.131749460043
$ws.Range("C26").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("I26").Value = 51
$ws.Range("K26").Value = 8.510638297872
$ws.Range("L26").Value = -5.555555555555
$ws.Range("C27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 70
$ws.Range("J27").Value = 73
$ws.Range("K27").Value = -4.109589041095
$ws.Range("L27").Value = 22.807017543859
$ws.Range("C28").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 37
$ws.Range("J28").Value = 54
$ws.Range("K28").Value = -31.481481481481
$ws.Range("L28").Value = -39.344262295082
$ws.Range("M28").Value = -7.5
$ws.Range("N28").Value = -71.969696969697
$ws.Range("C29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 27
$ws.Range("J29").Value = 45
$ws.Range("K29").Value = -40
$ws.Range("L29").Value = -50.909090909090
$ws.Range("M29").Value = -18.181818181818
$ws.Range("N29").Value = -75.892857142857
$ws.Range("G30").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
